# Scheduled-runner update: refresh market-price-derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several leve
# sheets. Only numeric value cells are touched; no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 232.77272
$ws.Range("I33").Value = 240.42857
$ws.Range("K33").Value = 240.42857
$ws.Range("M33").Value = -11.42857000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 89.166664
$ws.Range("I42").Value = 79.09090999999999
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 237.27273
$ws.Range("L42").Value = 600
$ws.Range("M42").Value = -7.272729999999967
$ws.Range("N42").Value = -1060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 55556180
$ws.Range("I80").Value = 401
$ws.Range("J80").Value = 200001200
$ws.Range("K80").Value = 1203
$ws.Range("L80").Value = 600003600
$ws.Range("M80").Value = -205
$ws.Range("N80").Value = -600005596

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 55556180
$ws.Range("I83").Value = 401
$ws.Range("J83").Value = 200001200
$ws.Range("K83").Value = 3609
$ws.Range("L83").Value = 1800010800
$ws.Range("M83").Value = 1383
$ws.Range("N83").Value = -1800020784

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 39698
$ws.Range("J87").Value = 39698
$ws.Range("L87").Value = 39698
$ws.Range("N87").Value = -42194

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 8536.462
$ws.Range("I88").Value = 890
$ws.Range("J88").Value = 10830.4
$ws.Range("K88").Value = 890
$ws.Range("L88").Value = 10830.4
$ws.Range("M88").Value = -484
$ws.Range("N88").Value = -11642.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 39698
$ws.Range("J90").Value = 39698
$ws.Range("L90").Value = 119094
$ws.Range("N90").Value = -131574

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 8536.462
$ws.Range("I91").Value = 890
$ws.Range("J91").Value = 10830.4
$ws.Range("K91").Value = 890
$ws.Range("L91").Value = 10830.4
$ws.Range("M91").Value = 514
$ws.Range("N91").Value = -13638.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 797.46155
$ws.Range("I125").Value = 360.66666
$ws.Range("J125").Value = 928.5
$ws.Range("K125").Value = 3245.99994
$ws.Range("L125").Value = 8356.5
$ws.Range("M125").Value = -785.9999399999997
$ws.Range("N125").Value = -13276.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2224586.5
$ws.Range("I137").Value = 4547435.5
$ws.Range("J137").Value = 2731.0435
$ws.Range("K137").Value = 13642306.5
$ws.Range("L137").Value = 8193.130500000001
$ws.Range("M137").Value = -13639756.5
$ws.Range("N137").Value = -13293.1305

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5267943
$ws.Range("I138").Value = 3735.9092
$ws.Range("J138").Value = 7412619.5
$ws.Range("K138").Value = 11207.7276
$ws.Range("L138").Value = 22237858.5
$ws.Range("M138").Value = -6067.7276
$ws.Range("N138").Value = -22248138.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 30000
$ws.Range("J56").Value = 30000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 11180421
$ws.Range("I74").Value = 16718626
$ws.Range("J74").Value = 104010.2
$ws.Range("K74").Value = 16718626
$ws.Range("L74").Value = 104010.2
$ws.Range("M74").Value = -16717752
$ws.Range("N74").Value = -105758.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 11180421
$ws.Range("I77").Value = 16718626
$ws.Range("J77").Value = 104010.2
$ws.Range("K77").Value = 83593130
$ws.Range("L77").Value = 520051
$ws.Range("M77").Value = -83588762
$ws.Range("N77").Value = -528787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8511790
$ws.Range("I132").Value = 13185865
$ws.Range("J132").Value = 53939.906
$ws.Range("K132").Value = 39557595
$ws.Range("L132").Value = 161819.718
$ws.Range("M132").Value = -39555065
$ws.Range("N132").Value = -166879.718

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37725.117
$ws.Range("I31").Value = 30474.656
$ws.Range("J31").Value = 48298.707
$ws.Range("K31").Value = 30474.656
$ws.Range("L31").Value = 48298.707
$ws.Range("M31").Value = -30179.656
$ws.Range("N31").Value = -48888.707

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 37725.117
$ws.Range("I34").Value = 30474.656
$ws.Range("J34").Value = 48298.707
$ws.Range("K34").Value = 30474.656
$ws.Range("L34").Value = 48298.707
$ws.Range("M34").Value = -30272.656
$ws.Range("N34").Value = -48702.707

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1407.6923
$ws.Range("I99").Value = 1375
$ws.Range("J99").Value = 1460
$ws.Range("K99").Value = 1375
$ws.Range("L99").Value = 1460
$ws.Range("M99").Value = 123
$ws.Range("N99").Value = -4456

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1407.6923
$ws.Range("I126").Value = 1375
$ws.Range("J126").Value = 1460
$ws.Range("K126").Value = 4125
$ws.Range("L126").Value = 4380
$ws.Range("M126").Value = -1655
$ws.Range("N126").Value = -9320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 100
$ws.Range("I48").Value = 100
$ws.Range("K48").Value = 300
$ws.Range("M48").Value = -50

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1819.1904
$ws.Range("I80").Value = 1375
$ws.Range("J80").Value = 1923.7059
$ws.Range("K80").Value = 4125
$ws.Range("L80").Value = 5771.1177
$ws.Range("M80").Value = -3189
$ws.Range("N80").Value = -7643.1177

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1819.1904
$ws.Range("I83").Value = 1375
$ws.Range("J83").Value = 1923.7059
$ws.Range("K83").Value = 12375
$ws.Range("L83").Value = 17313.3531
$ws.Range("M83").Value = -7695
$ws.Range("N83").Value = -26673.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2263.3555
$ws.Range("I132").Value = 1383.6111
$ws.Range("J132").Value = 2849.8518
$ws.Range("K132").Value = 12452.4999
$ws.Range("L132").Value = 25648.6662
$ws.Range("M132").Value = -9922.499900000001
$ws.Range("N132").Value = -30708.6662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2925.7827
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2925.7827
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 8777.348100000001
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -18977.3481

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27544.25
$ws.Range("J46").Value = 36665.555
$ws.Range("L46").Value = 36665.555
$ws.Range("N46").Value = -36977.555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 43317.42
$ws.Range("I70").Value = 63526.65
$ws.Range("J70").Value = 5144.4443
$ws.Range("K70").Value = 63526.65
$ws.Range("L70").Value = 5144.4443
$ws.Range("M70").Value = -63256.65
$ws.Range("N70").Value = -5684.4443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 43317.42
$ws.Range("I73").Value = 63526.65
$ws.Range("J73").Value = 5144.4443
$ws.Range("K73").Value = 63526.65
$ws.Range("L73").Value = 5144.4443
$ws.Range("M73").Value = -62590.65
$ws.Range("N73").Value = -7016.4443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4078.814
$ws.Range("J80").Value = 4138.718
$ws.Range("L80").Value = 4138.718
$ws.Range("N80").Value = -6134.718

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4078.814
$ws.Range("J83").Value = 4138.718
$ws.Range("L83").Value = 20693.59
$ws.Range("N83").Value = -30677.59

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2666.6667
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I46").Value = 10101574
$ws.Range("J46").Value = 600.6667
$ws.Range("K46").Value = 10101574
$ws.Range("L46").Value = 600.6667
$ws.Range("M46").Value = -10101386
$ws.Range("N46").Value = -976.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1792.3334
$ws.Range("I68").Value = 1699.8334
$ws.Range("K68").Value = 1699.8334
$ws.Range("M68").Value = -950.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1792.3334
$ws.Range("I71").Value = 1699.8334
$ws.Range("K71").Value = 8499.166999999999
$ws.Range("M71").Value = -4755.166999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1417.5714
$ws.Range("I93").Value = 1153.5
$ws.Range("J93").Value = 3002
$ws.Range("K93").Value = 1153.5
$ws.Range("L93").Value = 3002
$ws.Range("M93").Value = 94.5
$ws.Range("N93").Value = -5498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1685.7858
$ws.Range("I100").Value = 1422.4445
$ws.Range("J100").Value = 1810.5264
$ws.Range("K100").Value = 1422.4445
$ws.Range("L100").Value = 1810.5264
$ws.Range("M100").Value = -881.4445000000001
$ws.Range("N100").Value = -2892.5264

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 92882.48
$ws.Range("I136").Value = 55819.85
$ws.Range("J136").Value = 339966.66
$ws.Range("K136").Value = 167459.55
$ws.Range("L136").Value = 1019899.98
$ws.Range("M136").Value = -164909.55
$ws.Range("N136").Value = -1024999.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1411.125
$ws.Range("I126").Value = 1183.4286
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 3550.2858
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -1080.2858
$ws.Range("N126").Value = -13955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 57078.5
$ws.Range("J135").Value = 57078.5
$ws.Range("L135").Value = 57078.5
$ws.Range("N135").Value = -67218.5
